$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 13
$ws.Range("J2").Value = 13 / 14400
$ws.Range("K2").Value = 6396
$ws.Range("L2").Value = 6396 / 500000
